# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# This script:
#  1. Adds a new "Player Info" sheet before "ODI Batting" with player bio data.
#  2. Renames the MATCH_CARD_LINK column to MATCH_CODE on "ODI Batting" and
#     "ODI Bowling", replacing the full scorecard URL with just the numeric
#     match code.
#  3. Adds a new "ODI Batting Extra" sheet after "ODI Bowling" with additional
#     per-innings batting stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper style applier for header rows (bold, centered, top aligned, thin box)
# ---------------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------------
# 1. Insert the two new sheets first (note: inserting a sheet shifts the
#    positional identity of any previously-fetched sheet references, so all
#    Worksheets.Add() calls happen up-front, and every sheet we need to touch
#    afterwards is re-fetched fresh by name).
# ---------------------------------------------------------------------------
$battingSheetTmp = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetTmp)
$playerInfo.Name = "Player Info"

$bowlingSheetTmp = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Add($null, $bowlingSheetTmp)
$extraSheet.Name = "ODI Batting Extra"

# Re-fetch every sheet reference fresh, by name, now that the workbook has
# its final set of four sheets.
$playerInfo = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")

# ---------------------------------------------------------------------------
# 2. Fill in "Player Info".
# ---------------------------------------------------------------------------
$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    Set-HeaderStyle $cell
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4693"
$playerInfo.Cells.Item(2, 2).Value = "Thomas Kevin Curran"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------------
# 3. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, keep only the
#    numeric match code, and drop the stray empty INNING_NUMBER cells.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

for ($r = 2; $r -le 29; $r++) {
    $linkCell = $battingSheet.Cells.Item($r, 4)
    $linkVal = $linkCell.Value2
    if ($linkVal -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $matches[1]
    }

    $inningCell = $battingSheet.Cells.Item($r, 2)
    $inningVal = $inningCell.Value2
    if ($inningVal -eq $null -or $inningVal -eq "") {
        $inningCell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE, keep only the
#    numeric match code.
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le 29; $r++) {
    $linkCell = $bowlingSheet.Cells.Item($r, 2)
    $linkVal = $linkCell.Value2
    if ($linkVal -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $matches[1]
    }
}

# ---------------------------------------------------------------------------
# 5. Fill in "ODI Batting Extra".
# ---------------------------------------------------------------------------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $extraSheet.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    Set-HeaderStyle $cell
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4211", 9,      $null, $null, $null,     "NO"),
    @("4212", 9,      $null, $null, $null,     "NO"),
    @("4215", 10,     "0",   "0",   "0.76%",   "NO"),
    @("4254", 8,      "0",   "0",   $null,     "NO"),
    @("4260", 9,      "0",   "0",   $null,     "NO"),
    @("4284", 8,      "5",   "0",   "23.62%",  "NO"),
    @("4294", $null,  $null, $null, $null,     "NO"),
    @("4297", 8,      "5",   "0",   "9.09%",   "NO"),
    @("4300", 10,     "2",   "2",   "8.26%",   "NO"),
    @("4401", $null,  $null, $null, $null,     "NO"),
    @("4405", $null,  $null, $null, $null,     "NO"),
    @("4408", $null,  $null, $null, $null,     "NO"),
    @("4426", $null,  $null, $null, $null,     "NO"),
    @("4428", 9,      "4",   "0",   "11.59%",  "NO"),
    @("4430", 9,      "5",   "0",   "16.02%",  "NO"),
    @("4431", 8,      "0",   "1",   "6.29%",   "NO"),
    @("4454", 9,      "1",   "0",   "4.38%",   "NO"),
    @("4456", $null,  $null, $null, $null,     "NO"),
    @("4470", 8,      $null, $null, $null,     "NO"),
    @("4471", $null,  $null, $null, $null,     "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $codeCell = $extraSheet.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $posCell = $extraSheet.Cells.Item($r, 2)
    if ($row[1] -ne $null) {
        $posCell.Value = $row[1]
    }

    $num4Cell = $extraSheet.Cells.Item($r, 3)
    if ($row[2] -ne $null) {
        $num4Cell.NumberFormat = "@"
        $num4Cell.Value = $row[2]
    }

    $num6Cell = $extraSheet.Cells.Item($r, 4)
    if ($row[3] -ne $null) {
        $num6Cell.NumberFormat = "@"
        $num6Cell.Value = $row[3]
    }

    $pctCell = $extraSheet.Cells.Item($r, 5)
    if ($row[4] -ne $null) {
        $pctCell.NumberFormat = "@"
        $pctCell.Value = $row[4]
    }

    $momCell = $extraSheet.Cells.Item($r, 6)
    $momCell.Value = $row[5]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 6. Keep "Player Info" as the active sheet/tab.
# ---------------------------------------------------------------------------
$playerInfo.Activate()
